$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for the Price column (D) so numeric-looking
# strings like "0.999" or "6.70" are kept as text, matching the
# original inline-string cell type instead of being coerced to numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.824.33"
$ws.Range("E2").Value = "  +1.28%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.271.65"
$ws.Range("E3").Value = "  +0.82%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.17"
$ws.Range("E5").Value = "  +0.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "92.76"
$ws.Range("E6").Value = "  +1.23%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.531"
$ws.Range("E7").Value = "  +2.03%  "
$ws.Range("E8").Value = "  -0.15%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.486"
$ws.Range("E9").Value = "  +0.46%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.63"
$ws.Range("E10").Value = "  +2.01%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.66"
$ws.Range("E11").Value = "  -0.44%  "
$ws.Range("E12").Value = "  +0.48%  "
$ws.Range("E13").Value = "  -1.19%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.70"
$ws.Range("E14").Value = "  +1.44%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.623.78"
$ws.Range("E15").Value = "  +0.69%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.32"
$ws.Range("E16").Value = "  +1.34%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.275.66"
$ws.Range("E17").Value = "  +3.32%  "
$ws.Range("E18").Value = "  +3.94%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "41.738.61"
$ws.Range("E19").Value = "  +1.18%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.60"
$ws.Range("E20").Value = "  +3.57%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0908"
$ws.Range("E21").Value = "  +0.36%  "
$ws.Range("E22").Value = "  +0.90%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.14"
$ws.Range("E23").Value = "  +0.48%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "243.60"
$ws.Range("E24").Value = "  +1.56%  "
$ws.Range("E25").Value = "  +0.77%  "
$ws.Range("E26").Value = "  +4.02%  "
$ws.Range("E27").Value = "  +0.15%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "24.06"
$ws.Range("E28").Value = "  +1.47%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.53"
$ws.Range("E29").Value = "  -1.20%  "
$ws.Range("E30").Value = "  -5.36%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "35.37"
$ws.Range("E31").Value = "  +4.87%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "160.69"
$ws.Range("E32").Value = "  +1.46%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.25"
$ws.Range("E33").Value = "  +1.62%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0744"
$ws.Range("E35").Value = "  +1.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.02"
$ws.Range("E36").Value = "  -0.81%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "17.11"
$ws.Range("E37").Value = "  +3.36%  "
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.106"
$ws.Range("E38").Value = "  +2.47%  "
$ws.Range("B39").Value = "WEMIXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.36"
$ws.Range("E39").Value = "  -0.20%  "
$ws.Range("E41").Value = "  +1.68%  "
$ws.Range("E42").Value = "  -1.43%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.006.30"
$ws.Range("E43").Value = "  -2.78%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.54"
$ws.Range("E44").Value = "  -3.94%  "
$ws.Range("E45").Value = "  +2.30%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.30"
$ws.Range("E46").Value = "  +1.67%  "
$ws.Range("E47").Value = "  +2.07%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.91"
$ws.Range("E48").Value = "  -1.90%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "52.86"
$ws.Range("E49").Value = "  +3.54%  "
$ws.Range("E50").Value = "  +0.27%  "
$ws.Range("E51").Value = "  +1.08%  "
